$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by 23 rows (181 -> 204) for the new Paesi Bassi block
for ($i = 0; $i -lt 23; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Row 182
$ws.Range("A182").Value = "Paesi Bassi"
$ws.Range("C182").Value2 = 1
$ws.Range("D182").Value2 = 16939
$ws.Range("E182").Value2 = 17721
$ws.Range("F182").Value2 = 4
$ws.Range("G182").Value = "P"
$ws.Range("H182").Value = "Scioglimento anticipato per poter approvare riforma costituzionale postcoloniale"

# Row 183
$ws.Range("A183").Value = "Paesi Bassi"
$ws.Range("C183").Value2 = 2
$ws.Range("D183").Formula = "=E182"
$ws.Range("E183").Value2 = 19170
$ws.Range("F183").Value2 = 4
$ws.Range("G183").Value = "P"

# Row 184
$ws.Range("A184").Value = "Paesi Bassi"
$ws.Range("C184").Value2 = 3
$ws.Range("D184").Formula = "=E183"
$ws.Range("E184").Value2 = 20619
$ws.Range("F184").Value2 = 4
$ws.Range("G184").Value = "P"

# Row 185
$ws.Range("A185").Value = "Paesi Bassi"
$ws.Range("C185").Value2 = 4
$ws.Range("D185").Formula = "=E184"
$ws.Range("E185").Value2 = 21621
$ws.Range("F185").Value2 = 4
$ws.Range("G185").Value = "P"
$ws.Range("H185").Value = "Caduta del governo Drees III"

# Row 186
$ws.Range("A186").Value = "Paesi Bassi"
$ws.Range("C186").Formula = "=C185+1"
$ws.Range("D186").Formula = "=E185"
$ws.Range("E186").Value2 = 23146
$ws.Range("F186").Value2 = 4
$ws.Range("G186").Value = "P"

# Row 187
$ws.Range("A187").Value = "Paesi Bassi"
$ws.Range("C187").Formula = "=C186+1"
$ws.Range("D187").Formula = "=E186"
$ws.Range("E187").Value2 = 24518
$ws.Range("F187").Value2 = 4
$ws.Range("G187").Value = "P"
$ws.Range("H187").Value = "Caduta governo Cals"

# Row 188
$ws.Range("A188").Value = "Paesi Bassi"
$ws.Range("C188").Formula = "=C187+1"
$ws.Range("D188").Formula = "=E187"
$ws.Range("E188").Value2 = 26051
$ws.Range("F188").Value2 = 4
$ws.Range("G188").Value = "P"
$ws.Range("H188").Value = "Prime elezioni senza voto obbligatorio"

# Row 189
$ws.Range("A189").Value = "Paesi Bassi"
$ws.Range("C189").Formula = "=C188+1"
$ws.Range("D189").Formula = "=E188"
$ws.Range("E189").Value2 = 26632
$ws.Range("F189").Value2 = 4
$ws.Range("G189").Value = "P"
$ws.Range("H189").Value = "Caduta governo Biesheuvel I"

# Row 190
$ws.Range("A190").Value = "Paesi Bassi"
$ws.Range("C190").Formula = "=C189+1"
$ws.Range("D190").Formula = "=E189"
$ws.Range("E190").Value2 = 28270
$ws.Range("F190").Value2 = 5
$ws.Range("G190").Value = "P"
$ws.Range("H190").Value = "Mandato prolungato per scioglimento anticipato"

# Row 191
$ws.Range("A191").Value = "Paesi Bassi"
$ws.Range("C191").Formula = "=C190+1"
$ws.Range("D191").Formula = "=E190"
$ws.Range("E191").Value2 = 29732
$ws.Range("F191").Value2 = 4
$ws.Range("G191").Value = "P"

# Row 192
$ws.Range("A192").Value = "Paesi Bassi"
$ws.Range("C192").Formula = "=C191+1"
$ws.Range("D192").Formula = "=E191"
$ws.Range("E192").Value2 = 30202
$ws.Range("F192").Value2 = 4
$ws.Range("G192").Value = "P"
$ws.Range("H192").Value = "Caduta governo Van Agt II"

# Row 193
$ws.Range("A193").Value = "Paesi Bassi"
$ws.Range("C193").Formula = "=C192+1"
$ws.Range("D193").Formula = "=E192"
$ws.Range("E193").Value2 = 31553
$ws.Range("F193").Value2 = 4
$ws.Range("G193").Value = "P"

# Row 194
$ws.Range("A194").Value = "Paesi Bassi"
$ws.Range("C194").Formula = "=C193+1"
$ws.Range("D194").Formula = "=E193"
$ws.Range("E194").Value2 = 32757
$ws.Range("F194").Value2 = 4
$ws.Range("G194").Value = "P"

# Row 195
$ws.Range("A195").Value = "Paesi Bassi"
$ws.Range("C195").Formula = "=C194+1"
$ws.Range("D195").Formula = "=E194"
$ws.Range("E195").Value2 = 34457
$ws.Range("F195").Value2 = 5
$ws.Range("G195").Value = "P"
$ws.Range("H195").Value = "Mandato prolungato per scioglimento anticipato"

# Row 196
$ws.Range("A196").Value = "Paesi Bassi"
$ws.Range("C196").Formula = "=C195+1"
$ws.Range("D196").Formula = "=E195"
$ws.Range("E196").Value2 = 35921
$ws.Range("F196").Value2 = 4
$ws.Range("G196").Value = "P"

# Row 197
$ws.Range("A197").Value = "Paesi Bassi"
$ws.Range("C197").Formula = "=C196+1"
$ws.Range("D197").Formula = "=E196"
$ws.Range("E197").Value2 = 37391
$ws.Range("F197").Value2 = 4
$ws.Range("G197").Value = "P"
$ws.Range("H197").Value = "Caduta governo per rapporto su Srebrenica ma a camere quasi sciolte"

# Row 198
$ws.Range("A198").Value = "Paesi Bassi"
$ws.Range("C198").Formula = "=C197+1"
$ws.Range("D198").Formula = "=E197"
$ws.Range("E198").Value2 = 37643
$ws.Range("F198").Value2 = 4
$ws.Range("G198").Value = "P"
$ws.Range("H198").Value = "Caduta governo Balkenende I"

# Row 199
$ws.Range("A199").Value = "Paesi Bassi"
$ws.Range("C199").Formula = "=C198+1"
$ws.Range("D199").Formula = "=E198"
$ws.Range("E199").Value2 = 39043
$ws.Range("F199").Value2 = 4
$ws.Range("G199").Value = "P"
$ws.Range("H199").Value = "Caduta governo Balkenende II, non si erano estesi il mandato a 5 anni"

# Row 200
$ws.Range("A200").Value = "Paesi Bassi"
$ws.Range("C200").Formula = "=C199+1"
$ws.Range("D200").Formula = "=E199"
$ws.Range("E200").Value2 = 40338
$ws.Range("F200").Value2 = 5
$ws.Range("G200").Value = "P"
$ws.Range("H200").Value = "Caduta Balkenende IV, si erano estesi il mandato a 5"

# Row 201
$ws.Range("A201").Value = "Paesi Bassi"
$ws.Range("C201").Formula = "=C200+1"
$ws.Range("D201").Formula = "=E200"
$ws.Range("E201").Value2 = 41164
$ws.Range("F201").Value2 = 4
$ws.Range("G201").Value = "P"
$ws.Range("H201").Value = "Caduta governo Rutte I, non si erano estesi il mandato"

# Row 202
$ws.Range("A202").Value = "Paesi Bassi"
$ws.Range("C202").Formula = "=C201+1"
$ws.Range("D202").Formula = "=E201"
$ws.Range("E202").Value2 = 43034
$ws.Range("F202").Value2 = 5
$ws.Range("G202").Value = "P"
$ws.Range("H202").Value = "Mandato prolungato per scioglimento anticipato, Rutte II primo governo da Kok I a fare un interno mandato parlamentare"

# Row 203
$ws.Range("A203").Value = "Paesi Bassi"
$ws.Range("C203").Formula = "=C202+1"
$ws.Range("D203").Formula = "=E202"
$ws.Range("E203").Value2 = 44272
$ws.Range("F203").Value2 = 4
$ws.Range("G203").Value = "P"
$ws.Range("H203").Value = "Caduta governo Rutte III, elezioni non anticipate"

# Apply the date format used for the new date cells (including the trailing blank row 204)
$ws.Range("D182:E204").NumberFormat = "dd\-mmm\-yy"
